# "Adding feature files with cucumber"
#
# The test-data workbook's second sheet (the Cucumber feature-file /
# test-case matrix) had a typo fixed in one of its steps, the active
# selection moved, and the sheet's page setup (paper size / orientation)
# was configured.

$wb = $excel.ActiveWorkbook

# The sheet with the test-case table (dimension A2:E7) is the workbook's
# second sheet / tab.
$ws = $wb.Worksheets.Item(2)

# Fix the typo "prodcut" -> "product" in the scenario title held in B5.
$ws.Range("B5").Value = "As a non register user,I should NOT be able to refer a product to a friend"

# Move the active selection to D5.
$ws.Range("D5").Select()

# Configure the page setup for printing (paper size 9 = A4, portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
